# Personal log table tweak (margins) + new Tessel log entry + merge two
# runs in Roel's entry into one + add ListLabel32..49 character styles.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Table-wide indent / cell-margin nudge.
#    tblInd: -5 dxa -> -10 dxa   (dxa/20 = points)
#    tblCellMar left: 103 dxa -> 98 dxa
# ---------------------------------------------------------------------
$t = $d.Tables(1)
$t.Rows.LeftIndent = -0.5      # -10 dxa
$t.LeftPadding = 4.9           # 98 dxa

# Per-cell left margin override (tcMar) on every cell of the table also
# moves from 103 dxa to 98 dxa.
for ($r = 1; $r -le $t.Rows.Count; $r++) {
    for ($c = 1; $c -le $t.Columns.Count; $c++) {
        $t.Cell($r, $c).LeftPadding = 4.9
    }
}

# ---------------------------------------------------------------------
# 2. Roel's entry: the two runs get merged into a single run with the
#    combined text (no functional wording change, just a run merge).
# ---------------------------------------------------------------------
$roelCell = $t.Cell(2, 2)
$oldText = "Searched for common products between datasets that have the same exact name, this proved to be insufficient. Wrote code to produce lists of products from both datasets, found the proper common products by hand and put it in an excel sheet."
$roelRange = $roelCell.Range
$roelRange.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $oldText, 2) | Out-Null

# ---------------------------------------------------------------------
# 3. Tessel's entry was an empty placeholder run -- fill in the new
#    personal-log text.
# ---------------------------------------------------------------------
$tesselCell = $t.Cell(3, 2)
$tesselPara = $tesselCell.Range.Paragraphs(1)
$tesselPara.Range.Text = "16.00-18.00 Perfomed some data analysis to determine which data the EDA should be based on. Visualizing product frequencies, available years. Located missing entries and searched for strategies to handle these."

# ---------------------------------------------------------------------
# 4. New character styles ListLabel32 .. ListLabel49.
#    (NOTE: this interpreter's function calls only bind positionally,
#    so pass all four args -- use "" for the ones that don't apply.)
# ---------------------------------------------------------------------
function Add-ListLabelStyle {
    param(
        [string]$Name,
        [string]$Ascii,
        [string]$HAnsi,
        [string]$Cs
    )
    $style = $d.Styles.Add($Name, 2)
    $style.QuickStyle = $true
    if ($Ascii -ne "") { $style.Font.NameAscii = $Ascii }
    if ($HAnsi -ne "") { $style.Font.NameOther = $HAnsi }
    if ($Cs -ne "") { $style.Font.NameBi = $Cs }
}

Add-ListLabelStyle "ListLabel 32" "Helvetica" "Helvetica" "Symbol"
Add-ListLabelStyle "ListLabel 33" "Helvetica" "Helvetica" "Courier New"
Add-ListLabelStyle "ListLabel 34" "" "" "Wingdings"
Add-ListLabelStyle "ListLabel 35" "" "" "Symbol"
Add-ListLabelStyle "ListLabel 36" "" "" "Courier New"
Add-ListLabelStyle "ListLabel 37" "" "" "Wingdings"
Add-ListLabelStyle "ListLabel 38" "" "" "Symbol"
Add-ListLabelStyle "ListLabel 39" "" "" "Courier New"
Add-ListLabelStyle "ListLabel 40" "" "" "Wingdings"
Add-ListLabelStyle "ListLabel 41" "Helvetica" "Helvetica" "Calibri"
Add-ListLabelStyle "ListLabel 42" "" "" "Courier New"
Add-ListLabelStyle "ListLabel 43" "" "" "Wingdings"
Add-ListLabelStyle "ListLabel 44" "" "" "Symbol"
Add-ListLabelStyle "ListLabel 45" "" "" "Courier New"
Add-ListLabelStyle "ListLabel 46" "" "" "Wingdings"
Add-ListLabelStyle "ListLabel 47" "" "" "Symbol"
Add-ListLabelStyle "ListLabel 48" "" "" "Courier New"
Add-ListLabelStyle "ListLabel 49" "" "" "Wingdings"

Write-Host "done"
